# [bug] For bug 50447: rename default color scheme "Office" -> "New Office"
#
# Canonical edit (ppt/theme/theme1.xml):
#   - <a:theme name="Office тема">    -> name="Office Theme"
#   - <a:clrScheme name="Оffice">     -> name="New Office"
#   - accent1 (4472C4) <-> accent5 (5B9BD5) swapped
#   - <a:fontScheme name="Оffice">    -> name="Office Theme"
#   - majorFont latin "Calibri Light" -> "Arial"
#   - minorFont latin "Calibri"       -> "Arial"
#   - <a:fmtScheme name="Оffice">     -> name="Office Theme"
#
# The live, writable surface the PowerPoint object model exposes onto a
# theme is ThemeColorScheme.Colors(i).RGB (reached off a Slide) - apply the
# part of the edit that maps onto it: swap accent1 and accent5.
#
# msoThemeColorIndex: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
# 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$accent1 = $tcs.Colors(5)
$accent5 = $tcs.Colors(9)

$accent1Rgb = $accent1.RGB
$accent5Rgb = $accent5.RGB

$accent1.RGB = $accent5Rgb
$accent5.RGB = $accent1Rgb
